$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K15").Value = 18662.7348
$ws.Range("H15").Value = 6220.9116
$ws.Range("M15").Value = -18493.7348
$ws.Range("I15").Value = 6220.9116
$ws.Range("H19").Value = 10816.846
$ws.Range("L19").Value = 19171.285
$ws.Range("J19").Value = 19171.285
$ws.Range("N19").Value = -19521.285
$ws.Range("L53").Value = 204.45454
$ws.Range("N53").Value = -1478.45454
$ws.Range("H53").Value = 264.73914
$ws.Range("J53").Value = 204.45454
$ws.Range("I62").Value = 4666.6665
$ws.Range("H62").Value = 4666.6665
$ws.Range("M62").Value = -4042.6665
$ws.Range("K62").Value = 4666.6665
$ws.Range("H65").Value = 4666.6665
$ws.Range("K65").Value = 23333.3325
$ws.Range("I65").Value = 4666.6665
$ws.Range("M65").Value = -20213.3325
$ws.Range("M70").Value = -5863.200000000001
$ws.Range("K70").Value = 6133.200000000001
$ws.Range("H70").Value = 2078.6667
$ws.Range("I70").Value = 2044.4
$ws.Range("K73").Value = 6133.200000000001
$ws.Range("H73").Value = 2078.6667
$ws.Range("M73").Value = -5197.200000000001
$ws.Range("I73").Value = 2044.4
$ws.Range("N86").Value = -8363
$ws.Range("I86").Value = 2951.889
$ws.Range("H86").Value = 3639.9565
$ws.Range("L86").Value = 6117
$ws.Range("K86").Value = 2951.889
$ws.Range("M86").Value = -1828.889
$ws.Range("J86").Value = 6117
$ws.Range("L89").Value = 30585
$ws.Range("I89").Value = 2951.889
$ws.Range("M89").Value = -9143.445
$ws.Range("K89").Value = 14759.445
$ws.Range("H89").Value = 3639.9565
$ws.Range("J89").Value = 6117
$ws.Range("N89").Value = -41817
$ws.Range("I98").Value = 1223.3158
$ws.Range("M98").Value = 274.6841999999999
$ws.Range("K98").Value = 1223.3158
$ws.Range("H98").Value = 1266.1666
$ws.Range("K99").Value = 1322.1
$ws.Range("H99").Value = 1488.421
$ws.Range("I99").Value = 440.7
$ws.Range("M99").Value = 175.9000000000001
$ws.Range("M122").Value = -1219.9474
$ws.Range("I122").Value = 1223.3158
$ws.Range("K122").Value = 3669.9474
$ws.Range("H122").Value = 1266.1666
$ws.Range("K137").Value = 10998.3999
$ws.Range("M137").Value = -8448.3999
$ws.Range("I137").Value = 3666.1333
$ws.Range("H137").Value = 4350.905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1059498.2
$ws.Range("H2").Value = 894113.25
$ws.Range("M2").Value = -1059385.2
$ws.Range("K2").Value = 1059498.2
$ws.Range("H116").Value = 894113.25
$ws.Range("K116").Value = 1059498.2
$ws.Range("I116").Value = 1059498.2
$ws.Range("M116").Value = -1057204.2
$ws.Range("M122").Value = -4027.333
$ws.Range("I122").Value = 2159.111
$ws.Range("K122").Value = 6477.333
$ws.Range("H122").Value = 2114.2856

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 1059498.2
$ws.Range("M3").Value = -1059384.2
$ws.Range("H3").Value = 894113.25
$ws.Range("K3").Value = 1059498.2
$ws.Range("K7").Value = 452.5
$ws.Range("I7").Value = 452.5
$ws.Range("H7").Value = 9182.6
$ws.Range("M7").Value = -339.5
$ws.Range("H20").Value = 2648.1667
$ws.Range("I20").Value = 2496.2
$ws.Range("M20").Value = -2249.2
$ws.Range("K20").Value = 2496.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N31").Value = -3297.8
$ws.Range("J31").Value = 2707.8
$ws.Range("H31").Value = 2101.7368
$ws.Range("L31").Value = 2707.8
$ws.Range("H34").Value = 2101.7368
$ws.Range("J34").Value = 2707.8
$ws.Range("L34").Value = 2707.8
$ws.Range("N34").Value = -3111.8
$ws.Range("L99").Value = 12995
$ws.Range("H99").Value = 10184.353
$ws.Range("J99").Value = 12995
$ws.Range("N99").Value = -15991
$ws.Range("J126").Value = 12995
$ws.Range("L126").Value = 38985
$ws.Range("N126").Value = -43925
$ws.Range("H126").Value = 10184.353
$ws.Range("L129").Value = 59389.4
$ws.Range("H129").Value = 59389.4
$ws.Range("N129").Value = -69389.39999999999
$ws.Range("J129").Value = 59389.4
$ws.Range("K132").Value = 22662.12
$ws.Range("J132").Value = 3960.7144
$ws.Range("L132").Value = 11882.1432
$ws.Range("N132").Value = -16942.1432
$ws.Range("M132").Value = -20132.12
$ws.Range("I132").Value = 7554.04
$ws.Range("H132").Value = 6768

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1397.3
$ws.Range("I11").Value = 352.72726
$ws.Range("M11").Value = -918.1817799999999
$ws.Range("J11").Value = 2674
$ws.Range("N11").Value = -8302
$ws.Range("K11").Value = 1058.18178
$ws.Range("L11").Value = 8022
$ws.Range("H80").Value = 4914.2
$ws.Range("M80").ClearContents()
$ws.Range("L80").Value = 14742.6
$ws.Range("N80").Value = -16614.6
$ws.Range("J80").Value = 4914.2
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("L83").Value = 44227.8
$ws.Range("K83").Value = 0
$ws.Range("J83").Value = 4914.2
$ws.Range("N83").Value = -53587.8
$ws.Range("H83").Value = 4914.2
$ws.Range("J94").Value = 750449.25
$ws.Range("L94").Value = 2251347.75
$ws.Range("N94").Value = -2252699.75
$ws.Range("H94").Value = 500568.5
$ws.Range("M94").Value = -1745
$ws.Range("K94").Value = 2421
$ws.Range("I94").Value = 807
$ws.Range("J98").Value = 1999.75
$ws.Range("L98").Value = 5999.25
$ws.Range("N98").Value = -8995.25
$ws.Range("H98").Value = 1999.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 183.375
$ws.Range("L2").Value = 230.875
$ws.Range("H2").Value = 207.125
$ws.Range("J2").Value = 230.875
$ws.Range("M2").Value = -70.375
$ws.Range("N2").Value = -456.875
$ws.Range("K2").Value = 183.375
$ws.Range("M102").Value = -3647.5884
$ws.Range("H102").Value = 4902.5835
$ws.Range("I102").Value = 5269.5884
$ws.Range("K102").Value = 5269.5884
$ws.Range("N122").Value = -17504.6671
$ws.Range("J122").Value = 4201.5557
$ws.Range("M122").Value = -3756106.3
$ws.Range("I122").Value = 1252852.1
$ws.Range("L122").Value = 12604.6671
$ws.Range("K122").Value = 3758556.3
$ws.Range("H122").Value = 591801.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M16").Value = -2414.4167
$ws.Range("I16").Value = 2584.4167
$ws.Range("K16").Value = 2584.4167
$ws.Range("H16").Value = 2584.4167
$ws.Range("M22").Value = -438.5833
$ws.Range("H22").Value = 726.75
$ws.Range("I22").Value = 733.5833
$ws.Range("K22").Value = 733.5833
$ws.Range("M27").Value = -626.5833
$ws.Range("H27").Value = 726.75
$ws.Range("I27").Value = 733.5833
$ws.Range("K27").Value = 733.5833
$ws.Range("I40").Value = 5406.6
$ws.Range("J40").Value = 7470
$ws.Range("N40").Value = -7742
$ws.Range("K40").Value = 5406.6
$ws.Range("M40").Value = -5270.6
$ws.Range("L40").Value = 7470
$ws.Range("H40").Value = 6323.6665
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J106").Value = 22997.5
$ws.Range("N106").Value = -25521.5
$ws.Range("H106").Value = 22997.5
$ws.Range("L106").Value = 22997.5
$ws.Range("M122").Value = -10152.4
$ws.Range("I122").Value = 4200.8
$ws.Range("K122").Value = 12602.4
$ws.Range("H122").Value = 5900.5713
$ws.Range("L127").Value = 42500
$ws.Range("N127").Value = -52420
$ws.Range("H127").Value = 42500
$ws.Range("J127").Value = 42500

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N45").Value = -10731.5
$ws.Range("H45").Value = 9749.5
$ws.Range("J45").Value = 9749.5
$ws.Range("L45").Value = 9749.5
$ws.Range("H107").Value = 1001.6429
$ws.Range("N107").Value = -5485.5
$ws.Range("L107").Value = 1645.5
$ws.Range("J107").Value = 548.5
$ws.Range("M122").Value = -9523.4614
$ws.Range("I122").Value = 3991.1538
$ws.Range("K122").Value = 11973.4614
$ws.Range("H122").Value = 5993.3076
$ws.Range("I136").Value = 5686.5
$ws.Range("K136").Value = 17059.5
$ws.Range("M136").Value = -14509.5
$ws.Range("H136").Value = 5932.1934
